$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $Val) {
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = '58.305.78'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '2.995.16'
$ws.Range("E3").Value = '  +3.42%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws "D5" '559.45'
$ws.Range("E5").Value = '  +1.67%  '
Set-TextValue $ws "D6" '136.08'
$ws.Range("E6").Value = '  +11.86%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("E8").Value = '  +4.60%  '
$ws.Range("D9").Value = '2.983.71'
$ws.Range("E9").Value = '  +3.32%  '
Set-TextValue $ws "D10" '0.131'
$ws.Range("E10").Value = '  +4.69%  '
Set-TextValue $ws "D11" '4.87'
$ws.Range("E11").Value = '  +2.65%  '
Set-TextValue $ws "D12" '0.455'
$ws.Range("E12").Value = '  +5.06%  '
$ws.Range("E13").Value = '  +6.41%  '
Set-TextValue $ws "D14" '33.28'
$ws.Range("E14").Value = '  +5.27%  '
$ws.Range("E15").Value = '  +3.28%  '
$ws.Range("D16").Value = '3.487.28'
$ws.Range("E16").Value = '  +3.42%  '
$ws.Range("E17").Value = '  +9.36%  '
$ws.Range("D18").Value = '2.992.02'
$ws.Range("E18").Value = '  +3.52%  '
$ws.Range("D19").Value = '58.234.56'
$ws.Range("E19").Value = '  +1.04%  '
Set-TextValue $ws "D20" '423.58'
$ws.Range("E20").Value = '  +4.04%  '
Set-TextValue $ws "D21" '13.63'
$ws.Range("E21").Value = '  +6.63%  '
Set-TextValue $ws "D22" '0.712'
$ws.Range("E22").Value = '  +9.16%  '
$ws.Range("E23").Value = '  +7.19%  '
$ws.Range("E24").Value = '  +5.48%  '
Set-TextValue $ws "D25" '80.50'
$ws.Range("E25").Value = '  +4.83%  '
$ws.Range("E26").Value = '  -0.04%  '
Set-TextValue $ws "D27" '1.00'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("E28").Value = '  +2.60%  '
Set-TextValue $ws "D29" '2.07'
$ws.Range("E29").Value = '  +8.37%  '
$ws.Range("E30").Value = '  +6.21%  '
$ws.Range("E31").Value = '  +5.03%  '
Set-TextValue $ws "D32" '6.04'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("E33").Value = '  +3.04%  '
$ws.Range("E34").Value = '  +8.05%  '
$ws.Range("E35").Value = '  +7.14%  '
Set-TextValue $ws "D36" '2.09'
$ws.Range("E36").Value = '  +3.57%  '
$ws.Range("D37").Value = '0.0₃0721'
$ws.Range("E37").Value = '  +17.06%  '
$ws.Range("E38").Value = '  +6.10%  '
Set-TextValue $ws "D39" '48.51'
$ws.Range("E39").Value = '  +0.28%  '
Set-TextValue $ws "D40" '2.75'
$ws.Range("E40").Value = '  +17.32%  '
Set-TextValue $ws "D41" '398.02'
$ws.Range("E41").Value = '  +10.93%  '
Set-TextValue $ws "D42" '0.0351'
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("E43").Value = '  +2.56%  '
$ws.Range("D44").Value = '2.733.50'
$ws.Range("E44").Value = '  +4.75%  '
Set-TextValue $ws "D46" '125.25'
$ws.Range("E46").Value = '  +6.50%  '
Set-TextValue $ws "D47" '0.243'
$ws.Range("E47").Value = '  +6.91%  '
Set-TextValue $ws "D48" '2.02'
$ws.Range("E48").Value = '  +4.50%  '
$ws.Range("E49").Value = '  +2.99%  '
Set-TextValue $ws "D50" '23.29'
$ws.Range("E50").Value = '  +3.99%  '
$ws.Range("E51").Value = '  +4.32%  '
